$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 7292.1304
$ws.Range("I11").Value = 7292.1304
$ws.Range("K11").Value = 7292.1304
$ws.Range("M11").Value = -7152.1304

$ws.Range("H19").Value = 5542.2856
$ws.Range("I19").Value = 3624.5
$ws.Range("J19").Value = 8099.3335
$ws.Range("K19").Value = 3624.5
$ws.Range("L19").Value = 8099.3335
$ws.Range("M19").Value = -3449.5
$ws.Range("N19").Value = -8449.333500000001

$ws.Range("H38").Value = 353.125
$ws.Range("J38").Value = 991
$ws.Range("L38").Value = 2973
$ws.Range("N38").Value = -3717

$ws.Range("H40").Value = 2066.7778
$ws.Range("I40").Value = 1699.75
$ws.Range("J40").Value = 2360.4
$ws.Range("K40").Value = 1699.75
$ws.Range("L40").Value = 2360.4
$ws.Range("M40").Value = -1524.75
$ws.Range("N40").Value = -2710.4

$ws.Range("H50").Value = 22
$ws.Range("J50").Value = 22
$ws.Range("L50").Value = 66
$ws.Range("N50").Value = -1016

$ws.Range("H51").Value = 6275.6
$ws.Range("J51").Value = 5856.385
$ws.Range("L51").Value = 5856.385
$ws.Range("N51").Value = -6824.385

$ws.Range("H58").Value = 1509.8182
$ws.Range("J58").Value = 7505
$ws.Range("L58").Value = 22515
$ws.Range("N58").Value = -22815

$ws.Range("H70").Value = 602089.7
$ws.Range("J70").Value = 2184.182
$ws.Range("L70").Value = 6552.545999999999
$ws.Range("N70").Value = -7092.545999999999

$ws.Range("H73").Value = 602089.7
$ws.Range("J73").Value = 2184.182
$ws.Range("L73").Value = 6552.545999999999
$ws.Range("N73").Value = -8424.545999999998

$ws.Range("H80").Value = 618747.4
$ws.Range("I80").Value = 1039437.6
$ws.Range("J80").Value = 1735
$ws.Range("K80").Value = 3118312.8
$ws.Range("L80").Value = 5205
$ws.Range("M80").Value = -3117314.8
$ws.Range("N80").Value = -7201

$ws.Range("H83").Value = 618747.4
$ws.Range("I83").Value = 1039437.6
$ws.Range("J83").Value = 1735
$ws.Range("K83").Value = 9354938.4
$ws.Range("L83").Value = 15615
$ws.Range("M83").Value = -9349946.4
$ws.Range("N83").Value = -25599

$ws.Range("H112").Value = 78374.62
$ws.Range("I112").Value = 101378.9
$ws.Range("K112").Value = 304136.7
$ws.Range("M112").Value = -303028.7

$ws.Range("H132").Value = 3711.875
$ws.Range("I132").Value = 3601.5715
$ws.Range("J132").Value = 4484
$ws.Range("K132").Value = 10804.7145
$ws.Range("L132").Value = 13452
$ws.Range("M132").Value = -8274.7145
$ws.Range("N132").Value = -18512

$ws.Range("H135").Value = 58828236
$ws.Range("I135").Value = 66671630
$ws.Range("J135").Value = 2748
$ws.Range("K135").Value = 600044670
$ws.Range("L135").Value = 24732
$ws.Range("M135").Value = -600042135
$ws.Range("N135").Value = -29802

$ws.Range("H137").Value = 2440.9285
$ws.Range("I137").Value = 1886.4445
$ws.Range("J137").Value = 3439
$ws.Range("K137").Value = 5659.333500000001
$ws.Range("L137").Value = 10317
$ws.Range("M137").Value = -3109.333500000001
$ws.Range("N137").Value = -15417

$ws.Range("H138").Value = 2105.2432
$ws.Range("I138").Value = 1176.48
$ws.Range("K138").Value = 3529.44
$ws.Range("M138").Value = 1610.56

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2314.75
$ws.Range("I45").Value = 2074
$ws.Range("K45").Value = 2074
$ws.Range("M45").Value = -1697

$ws.Range("H61").Value = 27779764
$ws.Range("I61").Value = 33334164
$ws.Range("K61").Value = 33334164
$ws.Range("M61").Value = -33333952

$ws.Range("H88").Value = 9805657
$ws.Range("I88").Value = 23810332
$ws.Range("J88").Value = 2385
$ws.Range("K88").Value = 23810332
$ws.Range("L88").Value = 2385
$ws.Range("M88").Value = -23809926
$ws.Range("N88").Value = -3197

$ws.Range("H91").Value = 9805657
$ws.Range("I91").Value = 23810332
$ws.Range("J91").Value = 2385
$ws.Range("K91").Value = 23810332
$ws.Range("L91").Value = 2385
$ws.Range("M91").Value = -23808928
$ws.Range("N91").Value = -5193

$ws.Range("H136").Value = 27779764
$ws.Range("I136").Value = 33334164
$ws.Range("K136").Value = 100002492
$ws.Range("M136").Value = -99999942

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4178.8335
$ws.Range("I86").Value = 6413.4287
$ws.Range("K86").Value = 6413.4287
$ws.Range("M86").Value = -5290.4287

$ws.Range("H89").Value = 4178.8335
$ws.Range("I89").Value = 6413.4287
$ws.Range("K89").Value = 32067.1435
$ws.Range("M89").Value = -26451.1435

$ws.Range("H133").Value = 69000
$ws.Range("J133").Value = 69000
$ws.Range("L133").Value = 69000
$ws.Range("N133").Value = -79120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3165.625
$ws.Range("I31").Value = 1828.6428
$ws.Range("K31").Value = 1828.6428
$ws.Range("M31").Value = -1533.6428

$ws.Range("H34").Value = 3165.625
$ws.Range("I34").Value = 1828.6428
$ws.Range("K34").Value = 1828.6428
$ws.Range("M34").Value = -1626.6428

$ws.Range("H58").Value = 2005.2142
$ws.Range("I58").Value = 1809.3
$ws.Range("K58").Value = 1809.3
$ws.Range("M58").Value = -1606.3

$ws.Range("H134").Value = 2367.375
$ws.Range("I134").Value = 1647.6
$ws.Range("K134").Value = 4942.799999999999
$ws.Range("M134").Value = -2407.799999999999

$ws.Range("H136").Value = 2005.2142
$ws.Range("I136").Value = 1809.3
$ws.Range("K136").Value = 5427.9
$ws.Range("M136").Value = -2877.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 7151.5654
$ws.Range("J81").Value = 8445.789000000001
$ws.Range("L81").Value = 25337.367
$ws.Range("N81").Value = -27583.367

$ws.Range("H84").Value = 7151.5654
$ws.Range("J84").Value = 8445.789000000001
$ws.Range("L84").Value = 76012.10100000001
$ws.Range("N84").Value = -87244.10100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 57
$ws.Range("I2").Value = 44.42857
$ws.Range("K2").Value = 44.42857
$ws.Range("M2").Value = 68.57142999999999

$ws.Range("H49").Value = 41336
$ws.Range("J49").Value = 43603.2
$ws.Range("L49").Value = 43603.2
$ws.Range("N49").Value = -43971.2

$ws.Range("H99").Value = 6177.5
$ws.Range("I99").Value = 1413.2
$ws.Range("K99").Value = 1413.2
$ws.Range("M99").Value = 832.8

$ws.Range("H132").Value = 4425.5557
$ws.Range("I132").Value = 4047.762
$ws.Range("K132").Value = 12143.286
$ws.Range("M132").Value = -9613.286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1703.5385
$ws.Range("I46").Value = 921.1429000000001
$ws.Range("J46").Value = 2616.3333
$ws.Range("K46").Value = 921.1429000000001
$ws.Range("L46").Value = 2616.3333
$ws.Range("M46").Value = -733.1429000000001
$ws.Range("N46").Value = -2992.3333

$ws.Range("H55").Value = 485.55554
$ws.Range("I55").Value = 348.1
$ws.Range("J55").Value = 657.375
$ws.Range("K55").Value = 348.1
$ws.Range("L55").Value = 657.375
$ws.Range("M55").Value = -175.1
$ws.Range("N55").Value = -1003.375

$ws.Range("H74").Value = 55910.5
$ws.Range("I74").Value = 50197
$ws.Range("J74").Value = 57815
$ws.Range("K74").Value = 50197
$ws.Range("L74").Value = 57815
$ws.Range("M74").Value = -49199
$ws.Range("N74").Value = -59811

$ws.Range("H77").Value = 55910.5
$ws.Range("I77").Value = 50197
$ws.Range("J77").Value = 57815
$ws.Range("K77").Value = 150591
$ws.Range("L77").Value = 173445
$ws.Range("M77").Value = -145599
$ws.Range("N77").Value = -183429

$ws.Range("H136").Value = 4036.1538
$ws.Range("I136").Value = 3924.5715
$ws.Range("K136").Value = 11773.7145
$ws.Range("M136").Value = -9223.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2902
$ws.Range("J96").Value = 2999
$ws.Range("L96").Value = 2999
$ws.Range("N96").Value = -5745

$ws.Range("H132").Value = 3834.1177
$ws.Range("I132").Value = 3535.4285
$ws.Range("K132").Value = 10606.2855
$ws.Range("M132").Value = -8076.2855
